$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,6,7,8,10,11,12,13,15)

$data = @(
    @(12.41929701408701,11.8061794470366,9.551612687867879,35.52485888897037,37.96388835106357,16.76192063483097,11.29924557356578,8.440156849915784,11.29401446750008,14.95976690132833,26.70313410578052),
    @(12.19768129536238,11.82555937502406,9.542956447578382,35.63052537396351,38.10480249299925,16.81358878332315,11.32291074216963,8.258915295559136,11.29924546051333,14.9250641146095,26.79606946567962),
    @(12.06122574578846,11.83830731479546,9.538911335593102,35.70197492085089,38.20006924963298,16.84748764081838,11.33825053604416,8.146368358980958,11.30360991581294,14.90538018329118,26.8575762927958),
    @(12.00559296448437,11.84371611058421,9.537584119036705,35.73274203286574,38.24108570907687,16.86184907802713,11.34470566411022,8.10024625791926,11.30567913346021,14.89777233065736,26.88375828881805),
    @(11.99635568577965,11.84462717187356,9.537383187594914,35.73795055784903,38.24802886810021,16.86426686194136,11.34578987237016,8.092573898826764,11.3060403045729,14.89653417442879,26.88817328008736),
    @(12.06047547443749,11.8383793929071,9.538892133373929,35.7023831743584,38.20061353262007,16.84767910706334,11.33833676530347,8.14574730584328,11.30363664403246,14.90527590016157,26.85792486784155),
    @(12.34300500154923,11.81268583608115,9.548365629209005,35.55992841140585,38.01065788576125,16.77928487894737,11.3072376799638,8.377962073321115,11.29557940308034,14.94746757991589,26.73425580904906),
    @(12.89100351226886,11.76901074505003,9.576930977379652,35.33275551205404,37.70773763505554,16.6623895888251,11.25264930404137,8.820814127159784,11.28888704523896,15.04285023133193,26.52700535173678),
    @(13.2860523681435,11.74098152081356,9.60388151239747,35.19772769473185,37.52785971056686,16.58696993586295,11.21640864679067,9.13540035040643,11.28947087108632,15.12030303834328,26.39623853974871),
    @(13.46335488934754,11.72910498877301,9.617407271980992,35.1432335416404,37.45534753175409,16.55492372979495,11.20075390805709,9.275568227078509,11.29091896917137,15.15706361584318,26.34141850253214),
    @(13.53008809385313,11.72473283392569,9.622708374490495,35.12359545792196,37.42923258124442,16.54311346298461,11.19494486339422,9.32817737170776,11.29163625969655,15.17119658477405,26.32133076296243),
    @(13.51573496482913,11.72566889341192,9.621558764303256,35.12778048385317,37.43479706908568,16.54564257020477,11.19619065719511,9.316868636800415,11.29147428360575,15.16814345877882,26.32562714936576),
    @(13.46885357473425,11.72874278181084,9.617839822363477,35.14159790440105,37.45317208257546,16.55394558219963,11.20027361093032,9.279906076442158,11.29097460212409,15.1582221252332,26.33975241364496),
    @(13.44008254975929,11.73064192214968,9.615585112434189,35.15019142811213,37.46460244078227,16.5590737207815,11.20279002969943,9.257202981148502,11.29069049806148,15.15217248792182,26.34849199597985),
    @(13.2744114822154,11.74177522945454,9.603022804968761,35.20142857139169,37.53278637222771,16.58910972345573,11.21744840810675,9.126176837168366,11.28939993539158,15.11793077754398,26.39991509221083),
    @(13.17211650488892,11.74882869933357,9.59563847496125,35.23463680076632,37.57700399697549,16.60811503561403,11.22665342772812,9.045009477707213,11.28891030900999,15.0973107446101,26.43265696159376),
    @(13.11305602893578,11.7529679708287,9.59151052805475,35.25438959030326,37.60331335014632,16.61925937688574,11.23202618736598,8.99805011676319,11.28874006269091,15.08559480231028,26.45192842308818),
    @(13.09302276644869,11.7543836067508,9.590133454183555,35.26118954346847,37.61237165712036,16.62306925327077,11.2338587700366,8.982104896020473,11.28870158193763,15.08165295487214,26.45852883067896),
    @(13.18302952498633,11.74806933084036,9.596412219939079,35.23103421001976,37.57220620418334,16.60606984725868,11.2256654398962,9.053678593308931,11.28895091172742,15.09949091300574,26.42912607488124),
    @(13.48263530458028,11.72783651124933,9.618927325866457,35.13751230894722,37.4477383909737,16.551497970929,11.19907112083223,9.290775971797018,11.2911167952115,15.16113055037004,26.33558525590808),
    @(13.67604419550817,11.71534291834991,9.634685295570803,35.08220590559406,37.37422624000929,16.51772588335298,11.18238398071621,9.442977496666659,11.29351640722532,15.20265102303489,26.27836464572398),
    @(13.57305729702675,11.72194436656257,9.626180504863232,35.11119154033236,37.4127428074376,16.5355775360808,11.19122689653909,9.362011387435246,11.29214602804279,15.18038010556398,26.30854612784759),
    @(13.17809652176554,11.7484123792231,9.59606204393593,35.23266088261747,37.57437252042136,16.60699379803006,11.22611185787576,9.049760201932019,11.28893220863911,15.09850482603282,26.43072099416065),
    @(12.74381228344102,11.78011096531403,9.568146360771273,35.38861830196211,37.78220895503475,16.69217270055344,11.26673564051259,8.702687433664016,11.28972758570025,15.01572635981742,26.57929661637836)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowdata = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowdata[$j]
    }
}

Write-Host "done"